# Checkpoint 3 - country analysis
$wb = $excel.ActiveWorkbook

$wsTable2 = $wb.Worksheets.Item("Table - 2.1")
$wsTable3 = $wb.Worksheets.Item("Table-3.1")

# Fill in the new answer on "Table - 2.1" (post_ipo_debt? question)
$wsTable2.Range("C9").Value = "post_ipo_debt?"

# Fill in the Top/Second/Third English speaking country answers on "Table-3.1"
$wsTable3.Range("C5").Value = "United States of America"
$wsTable3.Range("C6").Value = "Canada"
$wsTable3.Range("C7").Value = "New Zealand"

# Move the active selection: Table - 2.1 no longer the active/selected sheet,
# selection on it resets to C10
$wsTable2.Range("C10").Select()

# Table-3.1 becomes the active (tab-selected) sheet with selection at C19
$wsTable3.Activate()
$wsTable3.Range("C19").Select()
